$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# fluid_mass (named range, Sheet1!$B$11) updated as part of integrating
# the Demand Model into ModelCenter. Q_cool (B12) recalculates from the
# existing formula automatically.
$ws.Range("B11").Value = 1.767144375
